$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Raise upper bound for Q3 (row 10) from 6.5 to 7.
# This replaces the shared formula in G8:G10 for G10 with a plain literal value.
$ws.Range("G10").Value = 7

# Lower the lower bound for Q6 (row 13) from 10 to 9.5.
$ws.Range("F13").Value = 9.5

# Update the sheet view: move the active selection to J20 (also clears the
# previously scrolled topLeftCell position).
$ws.Activate()
$ws.Range("J20").Select()

$wb.Application.CalculateFull()
